$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing rows down by 3 (insert 3 new rows at the top)
$ws.Range("A1:A3").EntireRow.Insert()

# New summary block in H1:J2
$ws.Range("H1").Value = "Total Person Hours Estimated:"
$ws.Range("I1").Value = 75

$ws.Range("H2").Value = "Actual Person Hours:"
$ws.Range("I2").Formula = "=SUM(C6:G13)"

$ws.Range("J1").Value = "We estimated 100 hours on the last project because we were building it from scratch. Here we estimated less because we're only adding a few features to an already complete project."

# Style (yellow fill + thin black border) for the two estimate cells -
# format I1 once, then copy/paste its format onto I2 so both cells end
# up sharing a single new cell-format entry (matches the saved file).
$ws.Range("I1").Interior.Color = 65535
$ws.Range("I1").Borders.LineStyle = 1
$ws.Range("I1").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New hours data points added to the existing grid
$ws.Range("F8").Value = 0.5
$ws.Range("F10").Value = 3.5
$ws.Range("F13").Value = 2

# Column widths for the new columns
$ws.Columns("H").ColumnWidth = 23.14
$ws.Columns("J").ColumnWidth = 139.43

# Restore selection to match the saved workbook state
$ws.Range("H13").Select()
